# Weekly update: a new price record for "Apio" (Feria Lagunitas de Puerto
# Montt) is inserted at row 174 (the sheet is sorted with the most recent
# entries near the top of this block). Inserting a whole row there pushes
# the existing rows 174-233 down to 175-234, carrying all their original
# values and formatting (incl. the date-formatted style on column D) with
# them - which is exactly what the target XML shows, right down to the
# former row 233 now living on at row 234.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("174").Insert()

$ws.Cells.Item(174, 1).Value  = 4
$ws.Cells.Item(174, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(174, 3).Value  = "Los Lagos"
$ws.Cells.Item(174, 4).Value  = 44663
$ws.Cells.Item(174, 5).Value  = 10
$ws.Cells.Item(174, 6).Value  = 100112017
$ws.Cells.Item(174, 7).Value  = "Apio"
$ws.Cells.Item(174, 8).Value  = "Americana (o)"
$ws.Cells.Item(174, 9).Value  = "Primera"
$ws.Cells.Item(174, 10).Value = 45
$ws.Cells.Item(174, 11).Value = 12000
$ws.Cells.Item(174, 12).Value = 12000
$ws.Cells.Item(174, 13).Value = 12000
$ws.Cells.Item(174, 14).Value = "`$/docena de matas"
$ws.Cells.Item(174, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(174, 16).Value = 2000
$ws.Cells.Item(174, 17).Value = 6
$ws.Cells.Item(174, 18).Value = "Hortaliza"
